# Round of edits/comments on manuscript:
#   Add Mike Ackerman's review comment anchored to "Figure 1" in the
#   first figure caption paragraph.

$d = $word.ActiveDocument

# Comments created through the Word object model are stamped with the
# active user's name/initials.
$word.UserName = "Mike Ackerman"
$word.UserInitials = "MA"

# Locate the literal text "Figure 1" (without the trailing ". ") inside
# the first figure caption so the comment range wraps only that text,
# just like the authored edit.
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("Figure 1", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $comment = $d.Comments.Add($rng, "It appears that the 2018 migration rate estimate is >100%. How would that be interpreted? Also, could the upper Cis could be capped at 100%?")
}

Write-Output "Comments: $($d.Comments.Count)"
